$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 0.0000000006650985193515199
$ws.Range("C9").Value = 0.000009435629509528259
$ws.Range("C10").Value = 0.00003940332640220782
$ws.Range("C11").Value = 0.00008912433899491021
$ws.Range("C12").Value = 0.000157591622060699
$ws.Range("C13").Value = 0.0002440099361728773
$ws.Range("C14").Value = 0.0003482974741566703
$ws.Range("C15").Value = 0.0004692287845529481
$ws.Range("D15").Value = 0.0000008081650430796593
$ws.Range("C16").Value = 0.0006076991141689865
$ws.Range("D16").Value = 0.00001229764095951468
$ws.Range("C17").Value = 0.0007643122737196944
$ws.Range("D17").Value = 0.00003568947197799329
$ws.Range("C18").Value = 0.0009379271276902408
$ws.Range("D18").Value = 0.00006928208659886968
$ws.Range("C19").Value = 0.001128317933738168
$ws.Range("D19").Value = 0.0001128035288148558
$ws.Range("C20").Value = 0.001336061316346792
$ws.Range("D20").Value = 0.000166174935932489
$ws.Range("C21").Value = 0.001561749592812818
$ws.Range("D21").Value = 0.0002295489303010083
$ws.Range("C22").Value = 0.001805448643154809
$ws.Range("D22").Value = 0.0003031009419451505
$ws.Range("C23").Value = 0.002067410665238578
$ws.Range("D23").Value = 0.0003867088377299903
$ws.Range("C24").Value = 0.002349753083193871
$ws.Range("D24").Value = 0.000480288412308698
$ws.Range("C25").Value = 0.002652335917633582
$ws.Range("D25").Value = 0.000584318451764875
$ws.Range("C26").Value = 0.002974199940300828
$ws.Range("D26").Value = 0.0006988194923462606
$ws.Range("C27").Value = 0.003317539269302034
$ws.Range("D27").Value = 0.0008239688789645609
$ws.Range("C28").Value = 0.003682649255309679
$ws.Range("D28").Value = 0.0009616568961458178
$ws.Range("C29").Value = 0.004069317874955292
$ws.Range("D29").Value = 0.00111192165371937
$ws.Range("C30").Value = 0.004478060312691045
$ws.Range("D30").Value = 0.001274711265054723
$ws.Range("C31").Value = 0.004910118576634573
$ws.Range("D31").Value = 0.001450009526464882
$ws.Range("C32").Value = 0.00536612163242128
$ws.Range("D32").Value = 0.001639223240251231
$ws.Range("C33").Value = 0.005846644121894883
$ws.Range("D33").Value = 0.001842673858989195
$ws.Range("C34").Value = 0.00635228118622505
$ws.Range("D34").Value = 0.002059382547000541
$ws.Range("C35").Value = 0.006883260272241341
$ws.Range("D35").Value = 0.002290251793841661
$ws.Range("C36").Value = 0.007440534474431162
$ws.Range("D36").Value = 0.002535418713155514
$ws.Range("C37").Value = 0.008022696785132575
$ws.Range("D37").Value = 0.002794240274016885
$ws.Range("C38").Value = 0.008631269815319626
$ws.Range("D38").Value = 0.003066559633381765
$ws.Range("C39").Value = 0.009268786461359937
$ws.Range("D39").Value = 0.003354782296716049
$ws.Range("C40").Value = 0.009934714595773051
$ws.Range("D40").Value = 0.003658734634046592
$ws.Range("C41").Value = 0.01062798312084516
$ws.Range("D41").Value = 0.003978072982383201
$ws.Range("C42").Value = 0.01135096559331785
$ws.Range("D42").Value = 0.004312605081372527
$ws.Range("C43").Value = 0.01210394534509394
$ws.Range("D43").Value = 0.004662640926836132
$ws.Range("C44").Value = 0.01288666543281509
$ws.Range("D44").Value = 0.005030447481218469
$ws.Range("C45").Value = 0.01370183163770503
$ws.Range("D45").Value = 0.005416823291020024
$ws.Range("C46").Value = 0.01455227264723075
$ws.Range("D46").Value = 0.005823459276959154
$ws.Range("C47").Value = 0.01543910533535238
$ws.Range("D47").Value = 0.006249238222572891
$ws.Range("C48").Value = 0.0163637961347358
$ws.Range("D48").Value = 0.006694855345445297
$ws.Range("C49").Value = 0.01732546117560708
$ws.Range("D49").Value = 0.007161776007343163
$ws.Range("C50").Value = 0.01832596889039937
$ws.Range("D50").Value = 0.007650794263614986
$ws.Range("C51").Value = 0.01936533390942373
$ws.Range("D51").Value = 0.008162314072178238
$ws.Range("C52").Value = 0.02044457714310993
$ws.Range("D52").Value = 0.008697393517711731
$ws.Range("C53").Value = 0.02156501356582434
$ws.Range("D53").Value = 0.009256626929320935
$ws.Range("C54").Value = 0.02272823272465302
$ws.Range("D54").Value = 0.009841366819386006
$ws.Range("C55").Value = 0.02393673450804536
$ws.Range("D55").Value = 0.01045132576196237
$ws.Range("C56").Value = 0.02519258430756122
$ws.Range("D56").Value = 0.01108843226762127
$ws.Range("C57").Value = 0.02649254998123601
$ws.Range("D57").Value = 0.01175148787787911
$ws.Range("C58").Value = 0.0278446210487839
$ws.Range("D58").Value = 0.01244329046937909
$ws.Range("C59").Value = 0.029250573213759
$ws.Range("D59").Value = 0.0131646659656236
$ws.Range("C60").Value = 0.03071670972925079
$ws.Range("D60").Value = 0.01391528352517645
$ws.Range("C61").Value = 0.03224358572193402
$ws.Range("D61").Value = 0.01469662418444839
$ws.Range("C62").Value = 0.03382887269645048
$ws.Range("D62").Value = 0.01551080629688504
$ws.Range("C63").Value = 0.03547629768586456
$ws.Range("D63").Value = 0.01635902259086846
$ws.Range("C64").Value = 0.0371953937326144
$ws.Range("D64").Value = 0.0172423184303316
$ws.Range("C65").Value = 0.03898588207687046
$ws.Range("D65").Value = 0.0181593763906759
$ws.Range("C66").Value = 0.04085762658167001
$ws.Range("D66").Value = 0.01911379600856392
$ws.Range("C67").Value = 0.04282157266545549
$ws.Range("D67").Value = 0.02010557070035519
$ws.Range("C68").Value = 0.04488444781750348
$ws.Range("D68").Value = 0.02113634446667834
$ws.Range("C69").Value = 0.04704517149529807
$ws.Range("D69").Value = 0.02220914553338983
$ws.Range("C70").Value = 0.04931876143688712
$ws.Range("D70").Value = 0.0233260620982971
$ws.Range("C71").Value = 0.05172621115150727
$ws.Range("D71").Value = 0.02448947349780582
$ws.Range("C72").Value = 0.05428015137236415
$ws.Range("D72").Value = 0.02570310134077462
$ws.Range("C73").Value = 0.05700474397298083
$ws.Range("D73").Value = 0.02697044173203458
$ws.Range("C74").Value = 0.05991986142938238
$ws.Range("D74").Value = 0.0282960432698273
$ws.Range("C75").Value = 0.06306903877812713
$ws.Range("D75").Value = 0.02968181932112346
$ws.Range("C76").Value = 0.06650529150148506
$ws.Range("D76").Value = 0.03113498291926217
$ws.Range("C77").Value = 0.07028224685155805
$ws.Range("D77").Value = 0.03266416662535409
$ws.Range("C78").Value = 0.07441681600291093
$ws.Range("D78").Value = 0.03428383540658114
$ws.Range("C79").Value = 0.07901137905182781
$ws.Range("D79").Value = 0.03600350373943971
$ws.Range("C80").Value = 0.08409467173774443
$ws.Range("D80").Value = 0.03783701239955043
$ws.Range("C81").Value = 0.08975035819458292
$ws.Range("D81").Value = 0.03979761293176122
$ws.Range("C82").Value = 0.09600584431973239
$ws.Range("D82").Value = 0.04189842413386412
$ws.Range("C83").Value = 0.1029425639265255
$ws.Range("D83").Value = 0.04418066703025838
$ws.Range("C84").Value = 0.1106348742668855
$ws.Range("D84").Value = 0.04668491221307414
$ws.Range("C85").Value = 0.1191382393332683
$ws.Range("D85").Value = 0.04944504913245163
$ws.Range("C86").Value = 0.1286349875012225
$ws.Range("D86").Value = 0.05256472747411913
$ws.Range("C87").Value = 0.1393352901057771
$ws.Range("D87").Value = 0.05618597459954322
$ws.Range("C88").Value = 0.1514286899542718
$ws.Range("D88").Value = 0.06069778075572794
$ws.Range("C89").Value = 0.1654801844953855
$ws.Range("D89").Value = 0.06812587048298349
$ws.Range("C90").Value = 0.1828980021416812
$ws.Range("D90").Value = 0.08328629220237169
$ws.Range("C91").Value = 0.2058925345994693
$ws.Range("D91").Value = 0.1076533263902451
$ws.Range("C92").Value = 0.2358592233270008
$ws.Range("D92").Value = 0.1403826682554797
$ws.Range("C93").Value = 0.2727950920125348
$ws.Range("D93").Value = 0.1809196242505155
$ws.Range("C94").Value = 0.317007260255766
$ws.Range("D94").Value = 0.2296609251037977
$ws.Range("C95").Value = 0.3682229550101819
$ws.Range("D95").Value = 0.2863549803326499
$ws.Range("C96").Value = 0.4269016881139613
$ws.Range("D96").Value = 0.35165210964099
$ws.Range("C97").Value = 0.4930472576756375
$ws.Range("D97").Value = 0.4254778085481455
$ws.Range("C98").Value = 0.5672890655118548
$ws.Range("D98").Value = 0.5087964645756605
$ws.Range("C99").Value = 0.6499683884910806
$ws.Range("D99").Value = 0.6018357027497487
$ws.Range("C100").Value = 0.7431053078412394
$ws.Range("D100").Value = 0.7069719452796176
$ws.Range("C101").Value = 0.8515803790060107
$ws.Range("D101").Value = 0.8303786856476056
